# Add par data for the new tournament "Linger Longer Invitational"
# (3 rounds x (18 holes + OUT + IN + RD) = 63 new rows appended after the
# existing data, which ends at row 1681).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$holeLabels = @("Hole 1","Hole 2","Hole 3","Hole 4","Hole 5","Hole 6","Hole 7","Hole 8","Hole 9","OUT","Hole 10","Hole 11","Hole 12","Hole 13","Hole 14","Hole 15","Hole 16","Hole 17","Hole 18","IN","RD")
$holeVals   = @(4,5,4,3,4,5,4,3,4,36,4,4,5,4,3,4,4,3,5,36,72)
$rounds     = @("Rd 1","Rd 2","Rd 3")
$tournament = "Linger Longer Invitational"

$row = 1682
for ($r = 0; $r -lt 3; $r++) {
    for ($i = 0; $i -lt 21; $i++) {
        $ws.Cells.Item($row, 1).Value = $tournament
        $ws.Cells.Item($row, 2).Value = $rounds[$r]
        $ws.Cells.Item($row, 3).Value = $holeLabels[$i]
        $ws.Cells.Item($row, 4).Value = $holeVals[$i]
        $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 3)).HorizontalAlignment = -4108
        $row = $row + 1
    }
}

# Stray centered/empty F:G cell formatting that trails along the new block
# (present on every new row except the very first and the 5th new row).
for ($row = 1683; $row -le 1745; $row++) {
    if ($row -ne 1686) {
        $ws.Cells.Item($row, 6).HorizontalAlignment = -4108
        $ws.Cells.Item($row, 7).HorizontalAlignment = -4108
    }
}

# Leave the view scrolled/selected near the newly entered data.
$ws.Range("G1736").Select()
